$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) data per the latest scrape.
# Rows 24 and 25 also swap ranking order (PEPE now above Litecoin).
# Price/volume columns are stored as plain text in this sheet, so force
# a Text number format before assigning values that look numeric (e.g.
# '165.60') to stop Excel from silently coercing them to real numbers
# and dropping the significant trailing zero.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.185.51'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.56%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.772.18'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.57%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '648.15'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.19%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.60'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.51%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.774.50'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.40%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.05%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.525'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.67%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.14%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.16%  '

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.74%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000239'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.33%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.95'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.02%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.401.34'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.67%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.765.76'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.02%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.042.59'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.78%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.80'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.85%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.01'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.14%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '467.01'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.11%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.58'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.61%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.707'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.48%  '

$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000144'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -6.24%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.81'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.23%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.30'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.26%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.20'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.17%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.10'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.55%  '

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.03%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.919.97'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.59%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.70'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.44%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.26'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.91%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.14'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.17%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.63'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.24%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.172'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +15.03%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.00%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.725.57'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.35%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.82'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.85%  '

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.12%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.77'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.60%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -7.52%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.16%  '

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.78%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '45.21'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.37%  '

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.25%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '155.82'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.75%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '47.30'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.20%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.296'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.81%  '

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.31%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.35'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.28%  '
